$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 held the text "R40" (a shared string). The commit updates it to
# the text value "1". A leading apostrophe forces Excel to store the
# numeric-looking entry as text (shared string) instead of a number,
# matching the original t="s" cell type.
$ws.Range("B11").Value = "'1"
